# Automatische test-sync: 2025-06-26 23:11:50
# Appends a new "Logs" row (row 33) for the automated test mail about
# opening hours, extends the conditional-formatting ranges that track the
# used range, and bumps the "Openingstijden / Locatie" tally on the
# "Dashboard" sheet from 9 to 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 33

$ws.Cells.Item($row, 1).Value = "Wanneer zijn jullie open?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #1: Wanneer zijn jullie open?"
$ws.Cells.Item($row, 4).Value = "Openingstijden / Locatie"
$ws.Cells.Item($row, 5).Value = "Beste klant,`r`n`r`nBedankt voor uw e-mail. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 tot 17:00 uur. We zijn gesloten in het weekend. Mocht u nog verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`r`n`r`nMet vriendelijke groet,`r`n[Naam bedrijf]"
$ws.Cells.Item($row, 6).Value = "2025-06-26 23:10:55"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"

# Setting multi-line text (E column) auto-pins an explicit row height;
# autofit it back to the sheet's implicit default so row 33 matches the
# unpinned rows above it (no explicit ht/customHeight attribute).
$ws.Rows.Item($row).AutoFit()

# Extend the conditional formatting sqref ranges (D/G/H/I) so they cover
# the newly added row, mirroring how Excel grows these ranges when a row
# is appended to the tracked table.
$ws.Range("D2:D32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D33"))
$ws.Range("G2:G32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G33"))
$ws.Range("H2:H32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H33"))
$ws.Range("I2:I32").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I33"))

# Dashboard summary: "Openingstijden / Locatie" count goes from 9 to 10.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 10
